# Add two new columns (I: "I0", J: "IF") with data for rows 2-60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (copy style from existing header cell H1 so formatting matches)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Row data: row, I-value, J-value
$data = @(
    @(2,6,7),
    @(3,5,5),
    @(4,7,7),
    @(5,8,8),
    @(6,6,6),
    @(7,10,10),
    @(8,8,8),
    @(9,8,8),
    @(10,5,5),
    @(11,8,8),
    @(12,6,6),
    @(13,7,7),
    @(14,7,7),
    @(15,10,10),
    @(16,4,5),
    @(17,7,7),
    @(18,8,8),
    @(19,6,6),
    @(20,7,7),
    @(21,6,6),
    @(22,7,8),
    @(23,8,8),
    @(24,6,6),
    @(25,8,8),
    @(26,9,9),
    @(27,9,9),
    @(28,5,6),
    @(29,9,9),
    @(30,7,7),
    @(31,6,6),
    @(32,9,9),
    @(33,9,9),
    @(34,7,7),
    @(35,6,6),
    @(36,6,7),
    @(37,9,9),
    @(38,8,8),
    @(39,6,6),
    @(40,5,5),
    @(41,7,7),
    @(42,8,8),
    @(43,5,6),
    @(44,6,6),
    @(45,7,7),
    @(46,8,8),
    @(47,6,6),
    @(48,6,6),
    @(49,7,7),
    @(50,6,7),
    @(51,9,9),
    @(52,8,8),
    @(53,8,8),
    @(54,9,9),
    @(55,9,9),
    @(56,7,7),
    @(57,9,9),
    @(58,9,9),
    @(59,6,6),
    @(60,8,8)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
